$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target final data (rows 2..15), columns A..J
# A=ID(num) B=Cislo objednavky C=Zakaznik D=Email E=Adresa F=Produkt
# G=Mnozstvi(num) H=Datum(text) I=Cislo faktury J=Stav
$data = @(
    @(14, "ORD014", "Hana VAVROVA",   "lukas.krumpach@gmail.com",  "Jaktáře 1475",                                  "🔥 HYALCHONDRO® HC PLUS", 1, "2025-03-25", "F014", "Nová"),
    @(13, "ORD013", "Лукас Крумпах",  "lukas.krumpach@gmail.com",  "Moscow, Потаповский переулок, д. 8/12 стр. 2", "🔥 HYALCHONDRO® HC PLUS", 1, "2025-03-25", "F013", "Zpracovává se"),
    @(12, "ORD012", "Mariia Isova",   "l.m.e.companycz@gmail.com", "Názovská 14",                                   "🔥 HYALCHONDRO® HC PLUS", 1, "2025-03-24", "F012", "Nová"),
    @(11, "ORD011", "Лукас Крумпах",  "lukas.krumpach@gmail.com",  "Moscow, Потаповский переулок, д. 8/12 стр. 2", "🔥 HYALCHONDRO® HC PLUS", 1, "2025-03-24", "F011", "Nová"),
    @(10, "ORD010", "Lukas Krumpach", "lukas.krumpach@gmail.com",  "Premonstrátů 910C",                             "🔥 HYALCHONDRO® HC PLUS", 1, "2025-03-24", "F010", "Nová"),
    @(9,  "ORD009", "Lukas Krumpach", "lukas.krumpach@gmail.com",  "Premonstrátů 910C",                             "🔥 HYALCHONDRO® HC PLUS", 1, "2025-03-24", "F009", "Nová"),
    @(8,  "ORD008", "Lukas Krumpach", "lukas.krumpach@gmail.com",  "Premonstrátů 910C",                             "🔥 HYALCHONDRO® HC PLUS", 1, "2025-03-23", "F008", "Nová"),
    @(7,  "ORD007", "Lukas Krumpach", "lukas.krumpach@gmail.com",  "Premonstrátů 910C",                             "🔥 HYALCHONDRO® HC PLUS", 1, "2025-03-23", "F007", "Nová"),
    @(6,  "ORD006", "Lukas Krumpach", "lukas.krumpach@gmail.com",  "Premonstrátů 910C",                             "🔥 HYALCHONDRO® HC PLUS", 1, "2025-03-23", "F006", "Nová"),
    @(5,  "ORD005", "Лукас Крумпах",  "lukas.krumpach@gmail.com",  "Moscow, Потаповский переулок, д. 8/12 стр. 2", "🔥 HYALCHONDRO® HC PLUS", 1, "2025-03-23", "F005", "Nová"),
    @(4,  "ORD004", "Lukas Krumpach", "lukas.krumpach@gmail.com",  "Premonstrátů 910C",                             "🔥 HYALCHONDRO® HC PLUS", 1, "2025-03-23", "F004", "Nová"),
    @(3,  "ORD003", "Lukas Krumpach", "lukas.krumpach@gmail.com",  "Premonstrátů 910C",                             "🔥 HYALCHONDRO® HC PLUS", 1, "2025-03-23", "F003", "Nová"),
    @(2,  "ORD002", "Лукас Крумпах",  "lukas.krumpach@gmail.com",  "Moscow, Потаповский переулок, д. 8/12 стр. 2", "🔥 HYALCHONDRO® HC PLUS", 1, "2025-03-23", "F002", "Nová"),
    @(1,  "ORD001", "Lukas Krumpach", "lukas.krumpach@gmail.com",  "Premonstrátů 910C",                             "🔥 HYALCHONDRO® HC PLUS", 1, "2025-03-23", "F001", "Nová")
)

$startRow = 2
$r = $startRow
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]

    # Date column stored as plain text in the source data (not a real date
    # serial) - force text format so Excel does not auto-convert it.
    $dateCell = $ws.Cells.Item($r, 8)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $row[7]

    $ws.Cells.Item($r, 9).Value = $row[8]
    $ws.Cells.Item($r, 10).Value = $row[9]

    $r++
}

$lastRow = $startRow + $data.Count - 1
$ws.Range("A1:J" + $lastRow).EntireColumn.AutoFit()
